$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New row 5/6: a "notes" label + the local-machine raw data row added above the chart header ---
$ws.Range("B5").Value = "5328217 ,418741 ,473308 ,461109 ,448229 ,448645 ,436339 ,432057"

$localData = @(5328217, 418741, 473308, 461109, 448229, 448645, 436339, 432057)
for ($i = 0; $i -lt $localData.Length; $i++) {
    $ws.Cells.Item(6, 2 + $i).Value = $localData[$i]
}

# --- Row 9 header text gains a "(1000)" qualifier ---
$ws.Range("A9").Value = "Ring Buffer - Size (1000)= Doesn't matter since addition/removal is only done at the top of the queue"

# --- Rows 13-18: relabel + fill in the measured series for each machine ---
$ws.Range("A13").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (Local) Locked"
$row13 = @(5328217, 418741, 473308, 461109, 448229, 448645, 436339, 432057)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $ws.Cells.Item(13, 2 + $i).Value = $row13[$i]
}

$ws.Range("A14").Value = "4 Core Intel Core i5-2500K CPU @ 3.30GHz (Local) C++ Spinlock"
$row14 = @(10604337, 10611323, 10564897, 10540146, 10570471, 10592364, 10539685, 10488351)
for ($i = 0; $i -lt $row14.Length; $i++) {
    $ws.Cells.Item(14, 2 + $i).Value = $row14[$i]
}

$ws.Range("A15").Value = "2 Core Intel Xeon CPU @ 3.00 GHz Locked (Ducss)"
$row15 = @(1321542, 1495327, 1467956, 1254631, 1363862, 1767556, 1358405, 1373794)
for ($i = 0; $i -lt $row15.Length; $i++) {
    $ws.Cells.Item(15, 2 + $i).Value = $row15[$i]
}

$ws.Range("A16").Value = "2 Core Intel Xeon CPU @ 3.00 GHz Spinlock (Ducss)"
$row16 = @(1191029, 1180264, 1123372, 1189681, 1200067, 1228259, 1205112, 1179244)
for ($i = 0; $i -lt $row16.Length; $i++) {
    $ws.Cells.Item(16, 2 + $i).Value = $row16[$i]
}

$ws.Range("A17").Value = "2 Core Intel Xeon CPU @ 2.80 GHz Locked (Netsoc)"
$row17 = @(3421422, 4856188, 4855373, 4557720, 4723686, 3093845, 3169190, 2640398)
for ($i = 0; $i -lt $row17.Length; $i++) {
    $ws.Cells.Item(17, 2 + $i).Value = $row17[$i]
}

$ws.Range("A18").Value = "2 Core Intel Xeon CPU @ 2.80 GHz Spinlock (Ducss)"
$row18 = @(2560144, 3529022, 2567171, 2586400, 2599223, 2575840, 2522506, 2445529)
for ($i = 0; $i -lt $row18.Length; $i++) {
    $ws.Cells.Item(18, 2 + $i).Value = $row18[$i]
}

# --- Update the sheet's active selection to match the new data extent ---
$excel.Application.Goto($ws.Range("A10:I18"))
